$wb = $excel.ActiveWorkbook

# ---- Step 1: insert new "2022-Q3" sheet before the current "2022-Q2" sheet ----
$wsOldQ2 = $wb.Worksheets.Item(2)
$wsNew = $wb.Worksheets.Add($wsOldQ2)
$wsNew.Name = "2022-Q3"

$ws1 = $wb.Worksheets.Item(1)   # 总计
$ws2 = $wb.Worksheets.Item(2)   # 2022-Q3 (the newly created sheet)

# ---- Step 2: header row (基金代码...仓位排名) ----
$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"

# copy the bold/border/center header style from sheet 总计 (B1) onto the new header row
$ws1.Range("B1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Step 3: fund holding rows (A=index, B=code, C=name, D..G text numbers, H=rank) ----
# row 2: 960010 工银核心价值混合H
$ws2.Range("A2").Value = 0
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "960010"
$ws2.Range("C2").Value = "工银核心价值混合H"
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "49.15"
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "85.16"
$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "3.69"
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = "1.8136"
$ws2.Range("H2").Value = 5

# row 3: 001008 工银国企改革主题股票
$ws2.Range("A3").Value = 1
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "001008"
$ws2.Range("C3").Value = "工银国企改革主题股票"
$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "8.66"
$ws2.Range("E3").NumberFormat = "@"
$ws2.Range("E3").Value = "81.53"
$ws2.Range("F3").NumberFormat = "@"
$ws2.Range("F3").Value = "2.67"
$ws2.Range("G3").NumberFormat = "@"
$ws2.Range("G3").Value = "0.2312"
$ws2.Range("H3").Value = 9

# row 4: 013049 兴业能源革新股票A
$ws2.Range("A4").Value = 2
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "013049"
$ws2.Range("C4").Value = "兴业能源革新股票A"
$ws2.Range("D4").NumberFormat = "@"
$ws2.Range("D4").Value = "4.63"
$ws2.Range("E4").NumberFormat = "@"
$ws2.Range("E4").Value = "88.84"
$ws2.Range("F4").NumberFormat = "@"
$ws2.Range("F4").Value = "4.10"
$ws2.Range("G4").NumberFormat = "@"
$ws2.Range("G4").Value = "0.1898"
$ws2.Range("H4").Value = 6

# row 5: 013050 兴业能源革新股票C
$ws2.Range("A5").Value = 3
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "013050"
$ws2.Range("C5").Value = "兴业能源革新股票C"
$ws2.Range("D5").NumberFormat = "@"
$ws2.Range("D5").Value = "3.65"
$ws2.Range("E5").NumberFormat = "@"
$ws2.Range("E5").Value = "88.84"
$ws2.Range("F5").NumberFormat = "@"
$ws2.Range("F5").Value = "4.10"
$ws2.Range("G5").NumberFormat = "@"
$ws2.Range("G5").Value = "0.1496"
$ws2.Range("H5").Value = 6

# row 6: 001672 国寿安保智慧生活股票
$ws2.Range("A6").Value = 4
$ws2.Range("B6").NumberFormat = "@"
$ws2.Range("B6").Value = "001672"
$ws2.Range("C6").Value = "国寿安保智慧生活股票"
$ws2.Range("D6").NumberFormat = "@"
$ws2.Range("D6").Value = "4.26"
$ws2.Range("E6").NumberFormat = "@"
$ws2.Range("E6").Value = "86.24"
$ws2.Range("F6").NumberFormat = "@"
$ws2.Range("F6").Value = "2.95"
$ws2.Range("G6").NumberFormat = "@"
$ws2.Range("G6").Value = "0.1257"
$ws2.Range("H6").Value = 7

# row 7: 010460 兴业研究精选混合A
$ws2.Range("A7").Value = 5
$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("B7").Value = "010460"
$ws2.Range("C7").Value = "兴业研究精选混合A"
$ws2.Range("D7").NumberFormat = "@"
$ws2.Range("D7").Value = "2.76"
$ws2.Range("E7").NumberFormat = "@"
$ws2.Range("E7").Value = "87.68"
$ws2.Range("F7").NumberFormat = "@"
$ws2.Range("F7").Value = "4.16"
$ws2.Range("G7").NumberFormat = "@"
$ws2.Range("G7").Value = "0.1148"
$ws2.Range("H7").Value = 6

# row 8: 200010 长城双动力混合A
$ws2.Range("A8").Value = 6
$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = "200010"
$ws2.Range("C8").Value = "长城双动力混合A"
$ws2.Range("D8").NumberFormat = "@"
$ws2.Range("D8").Value = "3.29"
$ws2.Range("E8").NumberFormat = "@"
$ws2.Range("E8").Value = "93.10"
$ws2.Range("F8").NumberFormat = "@"
$ws2.Range("F8").Value = "3.41"
$ws2.Range("G8").NumberFormat = "@"
$ws2.Range("G8").Value = "0.1122"
$ws2.Range("H8").Value = 4

# row 9: 910021 东方红启华三年持有期混合A
$ws2.Range("A9").Value = 7
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "910021"
$ws2.Range("C9").Value = "东方红启华三年持有期混合A"
$ws2.Range("D9").NumberFormat = "@"
$ws2.Range("D9").Value = "4.33"
$ws2.Range("E9").NumberFormat = "@"
$ws2.Range("E9").Value = "74.71"
$ws2.Range("F9").NumberFormat = "@"
$ws2.Range("F9").Value = "2.50"
$ws2.Range("G9").NumberFormat = "@"
$ws2.Range("G9").Value = "0.1082"
$ws2.Range("H9").Value = 10

# row 10: 015561 长城双动力混合C
$ws2.Range("A10").Value = 8
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "015561"
$ws2.Range("C10").Value = "长城双动力混合C"
$ws2.Range("D10").NumberFormat = "@"
$ws2.Range("D10").Value = "2.72"
$ws2.Range("E10").NumberFormat = "@"
$ws2.Range("E10").Value = "93.10"
$ws2.Range("F10").NumberFormat = "@"
$ws2.Range("F10").Value = "3.41"
$ws2.Range("G10").NumberFormat = "@"
$ws2.Range("G10").Value = "0.0928"
$ws2.Range("H10").Value = 4

# row 11: 004818 国寿安保目标策略灵活配置混合A
$ws2.Range("A11").Value = 9
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "004818"
$ws2.Range("C11").Value = "国寿安保目标策略灵活配置混合A"
$ws2.Range("D11").NumberFormat = "@"
$ws2.Range("D11").Value = "2.70"
$ws2.Range("E11").NumberFormat = "@"
$ws2.Range("E11").Value = "45.00"
$ws2.Range("F11").NumberFormat = "@"
$ws2.Range("F11").Value = "2.25"
$ws2.Range("G11").NumberFormat = "@"
$ws2.Range("G11").Value = "0.0608"
$ws2.Range("H11").Value = 4

# row 12: 004819 国寿安保目标策略灵活配置混合C
$ws2.Range("A12").Value = 10
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "004819"
$ws2.Range("C12").Value = "国寿安保目标策略灵活配置混合C"
$ws2.Range("D12").NumberFormat = "@"
$ws2.Range("D12").Value = "1.73"
$ws2.Range("E12").NumberFormat = "@"
$ws2.Range("E12").Value = "45.00"
$ws2.Range("F12").NumberFormat = "@"
$ws2.Range("F12").Value = "2.25"
$ws2.Range("G12").NumberFormat = "@"
$ws2.Range("G12").Value = "0.0389"
$ws2.Range("H12").Value = 4

# row 13: 002604 华夏新起点灵活配置混合A
$ws2.Range("A13").Value = 11
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "002604"
$ws2.Range("C13").Value = "华夏新起点灵活配置混合A"
$ws2.Range("D13").NumberFormat = "@"
$ws2.Range("D13").Value = "0.65"
$ws2.Range("E13").NumberFormat = "@"
$ws2.Range("E13").Value = "64.89"
$ws2.Range("F13").NumberFormat = "@"
$ws2.Range("F13").Value = "4.86"
$ws2.Range("G13").NumberFormat = "@"
$ws2.Range("G13").Value = "0.0316"
$ws2.Range("H13").Value = 6

# row 14: 011313 东方红启华三年持有期混合B
$ws2.Range("A14").Value = 12
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "011313"
$ws2.Range("C14").Value = "东方红启华三年持有期混合B"
$ws2.Range("D14").NumberFormat = "@"
$ws2.Range("D14").Value = "0.90"
$ws2.Range("E14").NumberFormat = "@"
$ws2.Range("E14").Value = "74.71"
$ws2.Range("F14").NumberFormat = "@"
$ws2.Range("F14").Value = "2.50"
$ws2.Range("G14").NumberFormat = "@"
$ws2.Range("G14").Value = "0.0225"
$ws2.Range("H14").Value = 10

# row 15: 002409 华夏新活力灵活配置混合A
$ws2.Range("A15").Value = 13
$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "002409"
$ws2.Range("C15").Value = "华夏新活力灵活配置混合A"
$ws2.Range("D15").NumberFormat = "@"
$ws2.Range("D15").Value = "0.15"
$ws2.Range("E15").NumberFormat = "@"
$ws2.Range("E15").Value = "69.89"
$ws2.Range("F15").NumberFormat = "@"
$ws2.Range("F15").Value = "4.01"
$ws2.Range("G15").NumberFormat = "@"
$ws2.Range("G15").Value = "0.0060"
$ws2.Range("H15").Value = 4

# row 16: 008213 华夏新起点灵活配置混合C
$ws2.Range("A16").Value = 14
$ws2.Range("B16").NumberFormat = "@"
$ws2.Range("B16").Value = "008213"
$ws2.Range("C16").Value = "华夏新起点灵活配置混合C"
$ws2.Range("D16").NumberFormat = "@"
$ws2.Range("D16").Value = "0.01"
$ws2.Range("E16").NumberFormat = "@"
$ws2.Range("E16").Value = "64.89"
$ws2.Range("F16").NumberFormat = "@"
$ws2.Range("F16").Value = "4.86"
$ws2.Range("G16").NumberFormat = "@"
$ws2.Range("G16").Value = "0.0005"
$ws2.Range("H16").Value = 6

# row 17: 015947 兴业研究精选混合C
$ws2.Range("A17").Value = 15
$ws2.Range("B17").NumberFormat = "@"
$ws2.Range("B17").Value = "015947"
$ws2.Range("C17").Value = "兴业研究精选混合C"
$ws2.Range("D17").NumberFormat = "@"
$ws2.Range("D17").Value = "0.00"
$ws2.Range("E17").NumberFormat = "@"
$ws2.Range("E17").Value = "87.68"
$ws2.Range("F17").NumberFormat = "@"
$ws2.Range("F17").Value = "4.16"
$ws2.Range("G17").Value = 0
$ws2.Range("H17").Value = 6

# row 18: 002410 华夏新活力灵活配置混合C
$ws2.Range("A18").Value = 16
$ws2.Range("B18").NumberFormat = "@"
$ws2.Range("B18").Value = "002410"
$ws2.Range("C18").Value = "华夏新活力灵活配置混合C"
$ws2.Range("D18").NumberFormat = "@"
$ws2.Range("D18").Value = "0.00"
$ws2.Range("E18").NumberFormat = "@"
$ws2.Range("E18").Value = "69.89"
$ws2.Range("F18").NumberFormat = "@"
$ws2.Range("F18").Value = "4.01"
$ws2.Range("G18").Value = 0
$ws2.Range("H18").Value = 4

# row 19: 481001 工银核心价值混合A
$ws2.Range("A19").Value = 17
$ws2.Range("B19").NumberFormat = "@"
$ws2.Range("B19").Value = "481001"
$ws2.Range("C19").Value = "工银核心价值混合A"
$ws2.Range("D19").NumberFormat = "@"
$ws2.Range("D19").Value = "-5.40"
$ws2.Range("E19").NumberFormat = "@"
$ws2.Range("E19").Value = "85.16"
$ws2.Range("F19").NumberFormat = "@"
$ws2.Range("F19").Value = "3.69"
$ws2.Range("G19").NumberFormat = "@"
$ws2.Range("G19").Value = "-0.1993"
$ws2.Range("H19").Value = 5

# ---- Step 4: update the 总计 (totals) summary sheet ----
# Existing rows 2-8 hold index 0-6; shift content down one quarter and
# append a new row 9 for the oldest quarter that falls off the bottom (2020-Q4, index 7).
$ws1.Range("A8").Copy()
$ws1.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 18
$ws1.Range("D2").Value = 2.9
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 11
$ws1.Range("D3").Value = 2.97
$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 13
$ws1.Range("D4").Value = 3.19
$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "2021-Q4"
$ws1.Range("C5").Value = 16
$ws1.Range("D5").Value = 4.97
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "2021-Q3"
$ws1.Range("C6").Value = 33
$ws1.Range("D6").Value = 7.93
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "2021-Q2"
$ws1.Range("C7").Value = 18
$ws1.Range("D7").Value = 7.72
$ws1.Range("A8").Value = 6
$ws1.Range("B8").Value = "2021-Q1"
$ws1.Range("C8").Value = 15
$ws1.Range("D8").Value = 3.53
$ws1.Range("A9").Value = 7
$ws1.Range("B9").Value = "2020-Q4"
$ws1.Range("C9").Value = 6
$ws1.Range("D9").Value = 1.76
